$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

$ws.Range("D2").Value = "28.797.72"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").Value = "1.881.10"
$ws.Range("E3").Value = "  +3.11%  "
$ws.Range("E4").Value = "  +0.52%  "
Set-TextValue $ws.Range("D5") "324.49"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("E6").Value = "  +0.47%  "
Set-TextValue $ws.Range("D7") "0.4676"
$ws.Range("E7").Value = "  +1.07%  "
Set-TextValue $ws.Range("D8") "0.3937"
$ws.Range("E8").Value = "  +2.37%  "
Set-TextValue $ws.Range("D9") "0.07934"
$ws.Range("E9").Value = "  +0.92%  "
Set-TextValue $ws.Range("D10") "0.9786"
$ws.Range("E10").Value = "  +1.99%  "
Set-TextValue $ws.Range("D11") "22.34"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").Value = "1.891.62"
$ws.Range("E12").Value = "  +5.26%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D13") "7.027"
$ws.Range("E13").Value = "  +2.55%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "5.745"
$ws.Range("E14").Value = "  +1.71%  "
Set-TextValue $ws.Range("D15") "0.06960"
$ws.Range("E15").Value = "  +1.38%  "
Set-TextValue $ws.Range("D16") "88.73"
$ws.Range("E16").Value = "  +2.66%  "
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "28.804.25"
$ws.Range("E21").Value = "  +2.88%  "
Set-TextValue $ws.Range("D22") "5.349"
$ws.Range("E22").Value = "  +0.85%  "
Set-TextValue $ws.Range("D23") "11.10"
$ws.Range("E23").Value = "  +1.19%  "
Set-TextValue $ws.Range("D24") "2.122"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").Value = "2.115.36"
$ws.Range("E25").Value = "  +4.67%  "
Set-TextValue $ws.Range("D26") "153.64"
$ws.Range("E26").Value = "  +0.97%  "
Set-TextValue $ws.Range("D27") "19.40"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("E28").Value = "  +0.19%  "
Set-TextValue $ws.Range("D29") "2.002"
$ws.Range("E29").Value = "  +1.58%  "
Set-TextValue $ws.Range("D30") "120.00"
$ws.Range("E30").Value = "  +2.78%  "
Set-TextValue $ws.Range("D31") "0.09401"
$ws.Range("E31").Value = "  +1.82%  "
Set-TextValue $ws.Range("D32") "0.9402"
$ws.Range("E32").Value = "  +0.54%  "
Set-TextValue $ws.Range("D33") "5.313"
$ws.Range("E33").Value = "  +0.69%  "
Set-TextValue $ws.Range("D34") "1.357"
$ws.Range("E34").Value = "  +3.25%  "
Set-TextValue $ws.Range("D35") "3.351"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("E36").Value = "  -0.13%  "
Set-TextValue $ws.Range("D37") "0.02122"
$ws.Range("E37").Value = "  -0.82%  "
Set-TextValue $ws.Range("D38") "1.159"
$ws.Range("E38").Value = "  +1.59%  "
Set-TextValue $ws.Range("D39") "7.903"
$ws.Range("E39").Value = "  +4.48%  "
Set-TextValue $ws.Range("D40") "0.5725"
$ws.Range("E40").Value = "  +2.61%  "
Set-TextValue $ws.Range("D42") "10.00"
$ws.Range("E42").Value = "  +0.79%  "
Set-TextValue $ws.Range("D43") "0.07328"
$ws.Range("E43").Value = "  +4.74%  "
Set-TextValue $ws.Range("D44") "11.80"
$ws.Range("E44").Value = "  +1.39%  "
Set-TextValue $ws.Range("D45") "0.5345"
$ws.Range("E45").Value = "  +2.03%  "
Set-TextValue $ws.Range("D46") "1.151"
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D47") "1.848"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D48") "2.108"
$ws.Range("E48").Value = "  -4.46%  "
Set-TextValue $ws.Range("D49") "114.16"
$ws.Range("E49").Value = "  +1.86%  "
Set-TextValue $ws.Range("D50") "2.368"
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("E51").Value = "  +0.63%  "
